$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H33").Value = 432.66666
$ws.Range("J33").Value = 650
$ws.Range("L33").Value = 650
$ws.Range("N33").Value = -1108

$ws.Range("H76").Value = 8622.23
$ws.Range("I76").Value = 10634.857
$ws.Range("J76").Value = 6274.1665
$ws.Range("K76").Value = 10634.857
$ws.Range("L76").Value = 6274.1665
$ws.Range("M76").Value = -10319.857
$ws.Range("N76").Value = -6904.1665

$ws.Range("H79").Value = 8622.23
$ws.Range("I79").Value = 10634.857
$ws.Range("J79").Value = 6274.1665
$ws.Range("K79").Value = 10634.857
$ws.Range("L79").Value = 6274.1665
$ws.Range("M79").Value = -9542.857
$ws.Range("N79").Value = -8458.166499999999

$ws.Range("H101").Value = 1025
$ws.Range("I101").Value = 765.625
$ws.Range("K101").Value = 2296.875
$ws.Range("M101").Value = -674.875

$ws.Range("H112").Value = 5557453.5
$ws.Range("J112").Value = 5557453.5
$ws.Range("L112").Value = 16672360.5
$ws.Range("N112").Value = -16674576.5

$ws.Range("H125").Value = 8000
$ws.Range("J125").Value = 8000
$ws.Range("L125").Value = 72000
$ws.Range("N125").Value = -76920

$ws.Range("H132").Value = 2681.9119
$ws.Range("I132").Value = 2245.0356
$ws.Range("K132").Value = 6735.1068
$ws.Range("M132").Value = -4205.1068

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("M34").ClearContents()

$ws.Range("H45").Value = 2093.2307
$ws.Range("I45").Value = 1913.4
$ws.Range("K45").Value = 1913.4
$ws.Range("M45").Value = -1536.4

$ws.Range("H74").Value = 71509500
$ws.Range("I74").Value = 83427576
$ws.Range("J74").Value = 1057
$ws.Range("K74").Value = 83427576
$ws.Range("L74").Value = 1057
$ws.Range("M74").Value = -83426702
$ws.Range("N74").Value = -2805

$ws.Range("H77").Value = 71509500
$ws.Range("I77").Value = 83427576
$ws.Range("J77").Value = 1057
$ws.Range("K77").Value = 417137880
$ws.Range("L77").Value = 5285
$ws.Range("M77").Value = -417133512
$ws.Range("N77").Value = -14021

$ws.Range("H110").Value = 22381.176
$ws.Range("I110").Value = 26356.428
$ws.Range("K110").Value = 26356.428
$ws.Range("M110").Value = -24311.428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2277.8
$ws.Range("I20").Value = 2445
$ws.Range("J20").Value = 2166.3333
$ws.Range("K20").Value = 2445
$ws.Range("L20").Value = 2166.3333
$ws.Range("M20").Value = -2198
$ws.Range("N20").Value = -2660.3333

$ws.Range("H31").Value = 500
$ws.Range("I31").Value = 500
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 500
$ws.Range("M31").Value = -248
$ws.Range("N31").Value = 0
$ws.Range("L31").ClearContents()

$ws.Range("H86").Value = 14049.087
$ws.Range("I86").Value = 5768.7617
$ws.Range("J86").Value = 100992.5
$ws.Range("K86").Value = 5768.7617
$ws.Range("L86").Value = 100992.5
$ws.Range("M86").Value = -4645.7617
$ws.Range("N86").Value = -103238.5

$ws.Range("H89").Value = 14049.087
$ws.Range("I89").Value = 5768.7617
$ws.Range("J89").Value = 100992.5
$ws.Range("K89").Value = 28843.8085
$ws.Range("L89").Value = 504962.5
$ws.Range("M89").Value = -23227.8085
$ws.Range("N89").Value = -516194.5

$ws.Range("H105").Value = 10959.5
$ws.Range("I105").Value = 13446
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 13446
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -11699
$ws.Range("N105").Value = -6994

$ws.Range("H141").Value = 42699.5
$ws.Range("J141").Value = 42699
$ws.Range("L141").Value = 42699
$ws.Range("N141").Value = -53059

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9797.532999999999
$ws.Range("I99").Value = 9961.857
$ws.Range("J99").Value = 7497
$ws.Range("K99").Value = 9961.857
$ws.Range("L99").Value = 7497
$ws.Range("M99").Value = -8463.857
$ws.Range("N99").Value = -10493

$ws.Range("H103").Value = 34250
$ws.Range("I103").Value = 8500
$ws.Range("K103").Value = 8500
$ws.Range("M103").Value = -7328

$ws.Range("H126").Value = 9797.532999999999
$ws.Range("I126").Value = 9961.857
$ws.Range("J126").Value = 7497
$ws.Range("K126").Value = 29885.571
$ws.Range("L126").Value = 22491
$ws.Range("M126").Value = -27415.571
$ws.Range("N126").Value = -27431

$ws.Range("H131").Value = 67796.164
$ws.Range("J131").Value = 67796.164
$ws.Range("L131").Value = 67796.164
$ws.Range("N131").Value = -77876.164

$ws.Range("H134").Value = 2998
$ws.Range("I134").Value = 2998
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 8994
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = -6459
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3943.5
$ws.Range("I94").Value = 2999
$ws.Range("J94").Value = 4888
$ws.Range("K94").Value = 8997
$ws.Range("L94").Value = 14664
$ws.Range("M94").Value = -8321
$ws.Range("N94").Value = -16016

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H132").Value = 3427.5356
$ws.Range("I132").Value = 3529.353
$ws.Range("J132").Value = 3270.182
$ws.Range("K132").Value = 10588.059
$ws.Range("L132").Value = 9810.545999999998
$ws.Range("M132").Value = -8058.059000000001
$ws.Range("N132").Value = -14870.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3717.6428
$ws.Range("I7").Value = 3540.6365
$ws.Range("K7").Value = 3540.6365
$ws.Range("M7").Value = -3428.6365

$ws.Range("H40").Value = 4182.9443
$ws.Range("I40").Value = 4080.875
$ws.Range("K40").Value = 4080.875
$ws.Range("M40").Value = -3944.875

$ws.Range("H61").Value = 2928.121
$ws.Range("I61").Value = 2335.0557
$ws.Range("K61").Value = 2335.0557
$ws.Range("M61").Value = -2133.0557

$ws.Range("H68").Value = 4320
$ws.Range("I68").Value = 2790
$ws.Range("K68").Value = 2790
$ws.Range("M68").Value = -2041

$ws.Range("H71").Value = 4320
$ws.Range("I71").Value = 2790
$ws.Range("K71").Value = 13950
$ws.Range("M71").Value = -10206

$ws.Range("H113").Value = 2928.121
$ws.Range("I113").Value = 2335.0557
$ws.Range("K113").Value = 2335.0557
$ws.Range("M113").Value = -165.0556999999999

$ws.Range("H122").Value = 4556.6665
$ws.Range("I122").Value = 3810.9
$ws.Range("J122").Value = 4995.353
$ws.Range("K122").Value = 11432.7
$ws.Range("L122").Value = 14986.059
$ws.Range("M122").Value = -8982.700000000001
$ws.Range("N122").Value = -19886.059

$ws.Range("H126").Value = 3717.6428
$ws.Range("I126").Value = 3540.6365
$ws.Range("K126").Value = 10621.9095
$ws.Range("M126").Value = -8151.9095

$ws.Range("H136").Value = 3034.9736
$ws.Range("I136").Value = 1974.6538
$ws.Range("J136").Value = 5332.3335
$ws.Range("K136").Value = 5923.9614
$ws.Range("L136").Value = 15997.0005
$ws.Range("M136").Value = -3373.9614
$ws.Range("N136").Value = -21097.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 16099
$ws.Range("I40").Value = 8000
$ws.Range("J40").Value = 21498.334
$ws.Range("K40").Value = 8000
$ws.Range("L40").Value = 21498.334
$ws.Range("M40").Value = -7851
$ws.Range("N40").Value = -21796.334

$ws.Range("H122").Value = 3403.3845
$ws.Range("I122").Value = 2925.1
$ws.Range("K122").Value = 8775.299999999999
$ws.Range("M122").Value = -6325.299999999999

$ws.Range("H126").Value = 19600
$ws.Range("I126").Value = 19600
$ws.Range("K126").Value = 58800
$ws.Range("M126").Value = -56330

$ws.Range("H136").Value = 1404
$ws.Range("I136").Value = 1088
$ws.Range("J136").Value = 3300
$ws.Range("K136").Value = 3264
$ws.Range("L136").Value = 9900
$ws.Range("M136").Value = -714
$ws.Range("N136").Value = -15000
